# Update column G ("K" - strikeouts) values on Sheet1 to match newly
# regenerated save_data (using K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 1
    6  = 0
    7  = 6
    8  = 1
    9  = 5
    10 = 0
    11 = 5
    12 = 5
    13 = 3
    14 = 7
    15 = 10
    16 = 6
    17 = 4
    18 = 3
    19 = 6
    20 = 3
    21 = 12
    22 = 3
    23 = 4
    24 = 6
    25 = 7
    26 = 2
    27 = 7
    28 = 0
    29 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
